# ============================================================================
# Atualiza a "Tabela de Presença nas Reuniões":
#   - Remove o Rodrigo da lista de nomes (linha inteira excluída)
#   - Adiciona mais 9 reuniões (12ª a 20ª), com as respectivas datas
#   - Marca presença (1) em todas as reuniões para os participantes restantes
#   - Ajusta zoom / seleção da planilha
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove a linha do Rodrigo (linha 6) - desloca Matheus/Miguel/Alysson
#    para cima automaticamente e remove "Rodrigo" da tabela de strings.
# ---------------------------------------------------------------------------
$ws.Rows(6).Delete()

# ---------------------------------------------------------------------------
# 2) Insere 9 novas colunas de reunião (N:V) antes da coluna "Total"
#    (a antiga coluna N "Total" passa a ser a coluna W).
# ---------------------------------------------------------------------------
$ws.Columns("N:V").Insert()

# ---------------------------------------------------------------------------
# 3) Cabeçalhos das novas reuniões (linha 3) e datas (linha 4)
# ---------------------------------------------------------------------------
$novosTitulos = @("12º Reunião","13º Reunião","14º Reunião","15º Reunião","16º Reunião","17º Reunião","18º Reunião","19º Reunião","20º Reunião")
$novasDatas   = @(44681,44688,44695,44702,44706,44707,44708,44709,44710)

$colunas = @("N","O","P","Q","R","S","T","U","V")
for ($i = 0; $i -lt $colunas.Length; $i++) {
    $col = $colunas[$i]
    # Os estilos (s=1 no cabeçalho, s=2 na data) já são herdados automaticamente
    # pela própria operação de Insert() das colunas, então basta escrever os valores.
    $ws.Range($col + "3").Value2 = $novosTitulos[$i]
    $ws.Range($col + "4").Value2 = $novasDatas[$i]
}

# ---------------------------------------------------------------------------
# 4) Marca presença "1" em todas as reuniões (C:V) para as 4 linhas restantes
#    (as reuniões 5 a 11, colunas G:M, que antes eram 0, agora também são 1)
# ---------------------------------------------------------------------------
for ($row = 5; $row -le 8; $row++) {
    for ($col = 7; $col -le 22; $col++) {   # G(7) .. V(22)
        $ws.Cells.Item($row, $col).Value2 = 1
    }
}

# ---------------------------------------------------------------------------
# 5) Ajusta zoom e seleção ativa, conforme o estado final salvo pelo autor
# ---------------------------------------------------------------------------
$ws.Range("E4").Select()
$excel.ActiveWindow.Zoom = 70

"OK"
